# WIP, non-working after adding IR calculations
#
# This script reorganizes the "Constants" block (Operating Altitude / Cruising
# Velocity move down one row to make room for the IR FOV input that used to
# share row 19 with the "IR Calculations" banner), fixes a copy/paste bug in
# the EO "Probability of Detecting a Vehicle" formulas (C6 -> C4), and adds a
# full parallel set of IR (infrared / night) detection calculations mirroring
# the existing EO (electro-optical / day) ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the EO "Probability of Detecting a Vehicle" block (rows 14-16): the
#    exponent term should divide by N50 (C4), not Dc_Vehicle (C6); and I14/I16
#    are simplified to reference the numerator cell directly.
# ---------------------------------------------------------------------------
$ws.Range("I14").Formula = "=1+I13"
$ws.Range("I15").Formula = "=(EO_N_Vehicle/0.75)^(2.7+0.7*(EO_N_Vehicle/C4))"
$ws.Range("I16").Formula = "=1+I15"

# ---------------------------------------------------------------------------
# 2. Move "Operating Altitude" (old B16/C16) down into row 17 (alongside the
#    existing "Time to scan a 5km by 5km box" formula that already lived
#    there), and move "Cruising Velocity" (old B17/C17) down into row 18.
#    Clear out the old B16/C16 cells since row 16 becomes EO-denominator-only.
# ---------------------------------------------------------------------------
$ws.Range("B16").ClearContents()
$c16Value = $ws.Range("C16").Value()
$ws.Range("C16").ClearContents()

$ws.Range("B17").Value = "Operating Altitude"
$ws.Range("C17").Value = $c16Value

$ws.Range("B18").Value = "Cuising Velocity"
$ws.Range("C18").Value = 100

# Repoint the defined names that track these two inputs at their new cells.
$wb.Names.Item("OperatingAltitude").RefersTo = "=Sheet1!`$C`$17"
$wb.Names.Item("CruisingVelocity").RefersTo = "=Sheet1!`$C`$18"

# ---------------------------------------------------------------------------
# 3. "IR Calculations" banner (F19) used to share row 19 with the "IR FOV"
#    input label (B19); that label now moves down to row 20 alongside a new
#    input value and the first IR GSD formula. Column B labels are entered
#    first (top to bottom), then column E/F, matching the original author's
#    shared-string insertion order.
# ---------------------------------------------------------------------------
$ws.Range("B19").ClearContents()

$ws.Range("B20").Value = "IR FOV"
$ws.Range("C20").Value = 45

$ws.Range("B21").Value = "IR Horizontal"
$ws.Range("C21").Value = 200

$ws.Range("B22").Value = "IR Vertical"
$ws.Range("C22").Value = 200

# New defined names for the IR inputs (needed before the GSD formulas below
# can reference IRFOV/IRHoriz/IRVert/OperatingAltitude).
$wb.Names.Add("IRFOV", "=Sheet1!`$C`$20")
$wb.Names.Add("IRHoriz", "=Sheet1!`$C`$21")
$wb.Names.Add("IRVert", "=Sheet1!`$C`$22")

$ws.Range("E20").Value = "IR GSDh"
$ws.Range("F20").Formula = "=2*TAN((IRFOV*PI()/180)/(2*IRHoriz))*OperatingAltitude*FeetToMeters"

$ws.Range("E21").Value = "IR GSDv"
$ws.Range("F21").Formula = "=2*TAN((IRFOV*PI()/180)/(2*IRVert))*OperatingAltitude*FeetToMeters"

$wb.Names.Add("IRGSDh", "=Sheet1!`$F`$20")
$wb.Names.Add("IRGSDv", "=Sheet1!`$F`$21")

# ---------------------------------------------------------------------------
# 4. IR_N_Human / IR_N_Vehicle (rows 23-24), mirroring EO_N_Human/EO_N_Vehicle.
# ---------------------------------------------------------------------------
$ws.Range("E23").Value = "IR_N_Human"
$ws.Range("F23").Formula = "=Dc_Human/(IRGSDh+IRGSDv)"

$ws.Range("E24").Value = "IR_N_Vehicle"
$ws.Range("F24").Formula = "=Dc_Vehicle/(IRGSDh+IRGSDv)"

$wb.Names.Add("IR_N_Human", "=Sheet1!`$F`$23")
$wb.Names.Add("IR_N_Vehicle", "=Sheet1!`$F`$24")

# ---------------------------------------------------------------------------
# 5. IR_GroundSwath / IR_Ground_Coverage_Rate (rows 26-27), mirroring
#    EO_GroundSwath / EO_Ground_Coverage_Rate.
# ---------------------------------------------------------------------------
$ws.Range("E26").Value = "IR_GroundSwath"
$ws.Range("F26").Formula = "=(TAN(0.5*IRFOV*PI()/180)-TAN(-0.5*IRFOV*PI()/180))*OperatingAltitude/3.281"

$ws.Range("E27").Value = "IR_Ground_Coverage_Rate"
$ws.Range("F27").Formula = "=IR_GroundSwath*CruisingVelocity*1.852/3.6"

$wb.Names.Add("IR_GroundSwath", "=Sheet1!`$F`$26")
$wb.Names.Add("IR_Ground_Coverage_Rate", "=Sheet1!`$F`$27")

# ---------------------------------------------------------------------------
# 6. Probability of Detecting a Human/Vehicle at Night (rows 29-32), mirroring
#    the EO daytime probability blocks (rows 13-16).
# ---------------------------------------------------------------------------
$ws.Range("E29").Value = "Probability of Detecting a Human at Night"
$ws.Range("H29").Value = "Numerator"
$ws.Range("I29").Formula = "=(IR_N_Human/0.75)^(2.7+0.7*(IR_N_Human/C4))"
$ws.Range("F29").Formula = "=I29/I30"

$ws.Range("H30").Value = "Denominator"
$ws.Range("I30").Formula = "=1+I29"

$ws.Range("E31").Value = "Probability of Detecting a Vehicle at Night"
$ws.Range("H31").Value = "Numerator"
$ws.Range("I31").Formula = "=(IR_N_Vehicle/0.75)^(2.7+0.7*(IR_N_Vehicle/C4))"
$ws.Range("F31").Formula = "=I31/I32"

$ws.Range("H32").Value = "Denominator"
$ws.Range("I32").Formula = "=1+I31"

$wb.Names.Add("Prob_Night_Human", "=Sheet1!`$F`$29")
$wb.Names.Add("Prob_Night_Vehicle", "=Sheet1!`$F`$31")

# ---------------------------------------------------------------------------
# 7. Time to scan a 5km by 5km box, split into explicit Day/Night names (the
#    existing T_Scan_5km_Box name is left alone, still pointing at F17), plus
#    the new night-time calculation in row 33.
# ---------------------------------------------------------------------------
$wb.Names.Add("T_Scan_5km_Box_Day", "=Sheet1!`$F`$17")

$ws.Range("E33").Value = "Time to scan a 5km by 5km box at Night"
$ws.Range("F33").Formula = "=25000000*60/IR_Ground_Coverage_Rate"

$wb.Names.Add("T_Scan_5km_Box_Night", "=Sheet1!`$F`$33")

# ---------------------------------------------------------------------------
# 8. Selection, as left by the author mid-edit.
# ---------------------------------------------------------------------------
$ws.Range("I10").Select()
